$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Uzair is More brilliant than Khizar"
$ws.Range("A2").Select()

$wb.Save()
